$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15").Value = "`nSat Jan 11 10:54:50.380 UTC`n---- node0_RP0_CPU0 ----`nMemory information:`n    Physical Memory     : 6144.0   MB`n    Free Memory         : 3548.386 MB`n    Memory State        :   Normal`n"

$ws.Range("C18").Value = "`nSat Jan 11 10:55:11.031 UTC`ncontroller Optics0/0/0/1`n description Scr1-Prot`n bits-per-symbol 3.4375`n rx-low-threshold -220`n logging events link-status`n!`n`n"

$ws.Range("C20").Value = "`nSat Jan 11 10:55:21.715 UTC"

$ws.Range("C23").Value = "`nSat Jan 11 10:55:42.390 UTC`nSSH version : Cisco-2.0 `n`nid       chan pty     location        state           userid    host                  ver authentication connection type`n-------------------------------------------------------------------------------------------------------------------------------`nIncoming sessions`n50479    1    vty0    0/RP0/CPU0      SESSION_OPEN    admin     10.110.195.179        v2  password       Command-Line-Interface `n`nOutgoing sessions`n`n"
